# Auto commit at 2025-11-13 8:05:27.77
# Updates the monthly/yearly/cumulative metric inputs on the "Metrics"
# sheet. The "today" sheet's B/E/F columns are live formulas that point
# back at these cells (=Metrics!B2 ... =Metrics!B13, plus the E/F running
# totals), so they recompute automatically once the source values change -
# no direct edits are needed there. Likewise A1's TODAY()-1 cache follows
# the host clock on recalc.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 160433.96000000002
$wsMetrics.Range("B3").Value  = 140832.94
$wsMetrics.Range("B4").Value  = 49908.89
$wsMetrics.Range("B5").Value  = 6635
$wsMetrics.Range("B6").Value  = 4956679.7100000009
$wsMetrics.Range("B7").Value  = 4182909.6200000006
$wsMetrics.Range("B8").Value  = 1456868.72
$wsMetrics.Range("B9").Value  = 192842
$wsMetrics.Range("B10").Value = 33422060.700000007
$wsMetrics.Range("B11").Value = 31458184.779999997
$wsMetrics.Range("B12").Value = 11738590.760000002
$wsMetrics.Range("B13").Value = 1290472

# Restore the recorded cursor position on the "Metrics" sheet.
$wsMetrics.Range("F27").Select()

# Re-select "today" last so it remains the active/visible tab, matching
# the saved workbook view, and park the cursor at its recorded cell.
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("E5").Select()
